{"js": "// Remove the internal spacing inside the `new_health_agent` merge field:\n//   {{ new_health_agent.name.full(middle=\u2019full\u2019) }}\n// becomes\n//   {{new_health_agent.name.full(middle=\u2019full\u2019)}}\n// (Leave the neighbouring `health_agent.name.full(...)` merge field, which\n// still has its surrounding spaces, untouched.)\n\nconst body = context.document.body;\n\nconst target = \"{{ new_health_agent.name.full(middle=\\u2019full\\u2019) }}\";\nconst replacement = \"{{new_health_agent.name.full(middle=\\u2019full\\u2019)}}\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the new_health_agent merge field to edit.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Remove the internal spacing inside the `new_health_agent` merge field:\n#   {{ new_health_agent.name.full(middle=\u2019full\u2019) }}\n# becomes\n#   {{new_health_agent.name.full(middle=\u2019full\u2019)}}\n# (The neighbouring `health_agent.name.full(...)` merge field keeps its\n# surrounding spaces and is left untouched.)\n\n$d = $word.ActiveDocument\n\n$rsquo = [char]0x2019\n\n$findText = \"{{ new_health_agent.name.full(middle=\" + $rsquo + \"full\" + $rsquo + \") }}\"\n$replaceText = \"{{new_health_agent.name.full(middle=\" + $rsquo + \"full\" + $rsquo + \")}}\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $findText\n$find.Replacement.Text = $replaceText\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
